# Insert a new data row at row 160 (pushing the existing rows 160-172 down
# to 161-173) and populate it with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 160..172 down to 161..173, leaving a blank row 160 to fill in.
$ws.Rows.Item(160).Insert()

$ws.Cells.Item(160, 1).Value = 10
$ws.Cells.Item(160, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(160, 3).Value = "La Araucanía"
$ws.Cells.Item(160, 4).Value = 44585
$ws.Cells.Item(160, 5).Value = 9
$ws.Cells.Item(160, 6).Value = 100112043
$ws.Cells.Item(160, 7).Value = "Pepino dulce"
$ws.Cells.Item(160, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(160, 9).Value = "Primera"
$ws.Cells.Item(160, 10).Value = 55
$ws.Cells.Item(160, 11).Value = 24000
$ws.Cells.Item(160, 12).Value = 24000
$ws.Cells.Item(160, 13).Value = 24000
$ws.Cells.Item(160, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(160, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(160, 16).Value = 1333
$ws.Cells.Item(160, 17).Value = 18
$ws.Cells.Item(160, 18).Value = "Hortaliza"
